$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Scratch cell used to force plain-number-looking strings to remain text,
# matching the source workbook where these are inline strings, not numbers.
$scratch = $ws.Range("ZZ1000")
$scratch.NumberFormat = "@"

$ws.Range("D2").Value = '38.715.24'
$ws.Range("E2").Value = '  +1.72%  '
$ws.Range("D3").Value = '2.090.73'
$ws.Range("E3").Value = '  -0.05%  '
$ws.Range("E4").Value = '  -0.03%  '
$scratch.Value = '229.08'
$scratch.Copy()
$ws.Range("D5").PasteSpecial(-4163)
$ws.Range("E5").Value = '  +0.08%  '
$scratch.Value = '0.616'
$scratch.Copy()
$ws.Range("D6").PasteSpecial(-4163)
$ws.Range("E6").Value = '  +0.08%  '
$scratch.Value = '61.13'
$scratch.Copy()
$ws.Range("D7").PasteSpecial(-4163)
$ws.Range("E7").Value = '  +0.46%  '
$ws.Range("E8").Value = '  -0.03%  '
$ws.Range("E9").Value = '  +1.14%  '
$scratch.Value = '0.0843'
$scratch.Copy()
$ws.Range("D10").PasteSpecial(-4163)
$ws.Range("E10").Value = '  +0.22%  '
$ws.Range("E11").Value = '  -0.15%  '
$scratch.Value = '15.19'
$scratch.Copy()
$ws.Range("D12").PasteSpecial(-4163)
$ws.Range("E12").Value = '  +4.16%  '
$ws.Range("D13").Value = '2.400.52'
$ws.Range("E13").Value = '  -0.07%  '
$scratch.Value = '21.99'
$scratch.Copy()
$ws.Range("D14").PasteSpecial(-4163)
$ws.Range("E14").Value = '  -0.47%  '
$scratch.Value = '0.814'
$scratch.Copy()
$ws.Range("D15").PasteSpecial(-4163)
$ws.Range("E15").Value = '  +5.17%  '
$ws.Range("E16").Value = '  -0.69%  '
$ws.Range("D17").Value = '2.083.87'
$ws.Range("E17").Value = '  -0.11%  '
$ws.Range("D18").Value = '38.657.78'
$ws.Range("E18").Value = '  +1.79%  '
$scratch.Value = '71.72'
$scratch.Copy()
$ws.Range("D19").PasteSpecial(-4163)
$ws.Range("E19").Value = '  +2.31%  '
$scratch.Value = '6.10'
$scratch.Copy()
$ws.Range("D20").PasteSpecial(-4163)
$ws.Range("E20").Value = '  +0.61%  '
$ws.Range("E21").Value = '  +0.13%  '
$scratch.Value = '227.43'
$scratch.Copy()
$ws.Range("D22").PasteSpecial(-4163)
$ws.Range("E22").Value = '  +1.50%  '
$ws.Range("E24").Value = '  -1.46%  '
$scratch.Value = '2.32'
$scratch.Copy()
$ws.Range("D25").PasteSpecial(-4163)
$ws.Range("E25").Value = '  +0.31%  '
$ws.Range("E26").Value = '  +1.63%  '
$scratch.Value = '170.97'
$scratch.Copy()
$ws.Range("D27").PasteSpecial(-4163)
$ws.Range("E27").Value = '  +0.80%  '
$scratch.Value = '0.140'
$scratch.Copy()
$ws.Range("D28").PasteSpecial(-4163)
$ws.Range("E28").Value = '  +5.82%  '
$ws.Range("E29").Value = '  +8.68%  '
$scratch.Value = '19.24'
$scratch.Copy()
$ws.Range("D30").PasteSpecial(-4163)
$ws.Range("E30").Value = '  +1.26%  '
$ws.Range("E31").Value = '  +5.22%  '
$ws.Range("E32").Value = '  +0.25%  '
$ws.Range("E33").Value = '  +1.71%  '
$ws.Range("E34").Value = '  +0.70%  '
$scratch.Value = '0.0609'
$scratch.Copy()
$ws.Range("D35").PasteSpecial(-4163)
$ws.Range("E35").Value = '  +0.30%  '
$scratch.Value = '6.54'
$scratch.Copy()
$ws.Range("D36").PasteSpecial(-4163)
$ws.Range("E36").Value = '  +0.87%  '
$ws.Range("E37").Value = '  -0.54%  '
$scratch.Value = '3.57'
$scratch.Copy()
$ws.Range("D38").PasteSpecial(-4163)
$ws.Range("E38").Value = '  +1.40%  '
$ws.Range("E39").Value = '  -0.07%  '
$scratch.Value = '17.99'
$scratch.Copy()
$ws.Range("D40").PasteSpecial(-4163)
$ws.Range("E40").Value = '  -0.14%  '
$ws.Range("E41").Value = '  +4.76%  '
$ws.Range("B42").Value = 'Maker'
$ws.Range("C42").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D42").Value = '1.533.76'
$ws.Range("E42").Value = '  -0.72%  '
$ws.Range("B43").Value = 'Aave'
$ws.Range("C43").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$scratch.Value = '100.82'
$scratch.Copy()
$ws.Range("D43").PasteSpecial(-4163)
$ws.Range("E43").Value = '  +0.46%  '
$ws.Range("E44").Value = '  -0.79%  '
$ws.Range("E45").Value = '  +1.30%  '
$scratch.Value = '7.72'
$scratch.Copy()
$ws.Range("D46").PasteSpecial(-4163)
$ws.Range("E46").Value = '  +7.10%  '
$ws.Range("E47").Value = '  +0.95%  '
$scratch.Value = '4.08'
$scratch.Copy()
$ws.Range("D48").PasteSpecial(-4163)
$ws.Range("E48").Value = '  -1.36%  '
$ws.Range("E49").Value = '  +0.94%  '
$ws.Range("E50").Value = '  -1.00%  '
$ws.Range("D51").Value = '2.286.66'
$ws.Range("E51").Value = '  -0.06%  '

$scratch.Clear()
$excel.Application.CutCopyMode = $false
